# PROS-7403 Sanofi Brand Block removed from all templates
#
# The "Primary_Brand_Blocking" sheet contained 10 rows (3-12) of KPI data
# for the "Blocked Together" / "Blocked Together Per Brand" KPI tied to the
# Sanofi brand block, which has been removed from the template. Deleting the
# rows shifts all subsequent rows (former 13-47) up to become rows 3-37.

$wb = $excel.ActiveWorkbook

$wsBlocking = $wb.Worksheets.Item("Primary_Brand_Blocking")

# Remove the 10 "Blocked Together" rows (rows 3 through 12 inclusive).
$wsBlocking.Rows("3:12").Delete()

# The workbook now opens focused on the Primary_Brand_Blocking sheet instead
# of KPIs, with the selection resting on A3 (the first surviving data row).
$wsBlocking.Select()
$wsBlocking.Range("A3").Select()
